$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.783.62"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "2.479.64"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'319.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").Value = "'93.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("E7").Value = "  +1.03%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.520"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.0869"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.03%  "
$ws.Range("B11").Value = "Avalanche"
$ws.Range("C11").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D11").Value = "'33.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.73%  "
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Value = "2.862.03"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "'6.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "'15.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").Value = "2.476.90"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "'0.795"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.30%  "
$ws.Range("D18").Value = "41.751.62"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").Value = "'6.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").Value = "0.0₃0953"
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("D21").Value = "'71.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "'11.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.50%  "
$ws.Range("D23").Value = "'240.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("E25").Value = "  +2.70%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "'24.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("D29").Value = "'9.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.77%  "
$ws.Range("D30").Value = "'36.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.16%  "
$ws.Range("D31").Value = "'158.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("D36").Value = "'17.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("D37").Value = "'1.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.08%  "
$ws.Range("D38").Value = "'2.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.54%  "
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("E42").Value = "  +10.39%  "
$ws.Range("D43").Value = "1.992.20"
$ws.Range("E43").Value = "  +2.68%  "
$ws.Range("D44").Value = "'19.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.76%  "
$ws.Range("D45").Value = "'0.0287"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("D46").Value = "'3.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.75%  "
$ws.Range("E47").Value = "  +3.73%  "
$ws.Range("D48").Value = "2.718.24"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("D49").Value = "'97.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").Value = "'74.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.93%  "
$ws.Range("D51").Value = "'67.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.87%  "
